$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-10-07T18:07:14"

$ws.Range("V4").Value = 97.73999999999999
$ws.Range("W4").Value = 70.27
$ws.Range("X4").Value = 46.5
$ws.Range("Y4").Value = 57.46
$ws.Range("V6").Value = -2.35
$ws.Range("W6").Value = -1.48
$ws.Range("X6").Value = -0.5600000000000001
$ws.Range("Y6").Value = -0.86
$ws.Range("V9").Value = 102.66
$ws.Range("W9").Value = 73.20999999999999
$ws.Range("X9").Value = 48.92
$ws.Range("Y9").Value = 60.19
$ws.Range("V11").Value = 2.57
$ws.Range("W11").Value = 1.46
$ws.Range("X11").Value = 1.86
$ws.Range("Y11").Value = 1.87
$ws.Range("V14").Value = 116.62
$ws.Range("W14").Value = 153.24
$ws.Range("X14").Value = 176.47
$ws.Range("Y14").Value = 60.19
$ws.Range("Z14").Value = 37.15
$ws.Range("V15").Value = 13.86
$ws.Range("W15").Value = 80.03
$ws.Range("X15").Value = 127.49
$ws.Range("V16").Value = 2.67
$ws.Range("W16").Value = 1.46
$ws.Range("X16").Value = 1.91
$ws.Range("Y16").Value = 1.87
$ws.Range("Z16").Value = 1.52
$ws.Range("V19").Value = 97.55
$ws.Range("W19").Value = 70.2
$ws.Range("X19").Value = 46.46
$ws.Range("Y19").Value = 57.4
$ws.Range("V21").Value = -2.54
$ws.Range("W21").Value = -1.54
$ws.Range("X21").Value = -0.6
$ws.Range("Y21").Value = -0.92
$ws.Range("V24").Value = 97.55
$ws.Range("W24").Value = 70.2
$ws.Range("X24").Value = 46.46
$ws.Range("Y24").Value = 57.4
$ws.Range("V26").Value = -2.54
$ws.Range("W26").Value = -1.54
$ws.Range("X26").Value = -0.6
$ws.Range("Y26").Value = -0.92
$ws.Range("V29").Value = 96.98999999999999
$ws.Range("W29").Value = 70
$ws.Range("X29").Value = 46.32
$ws.Range("Y29").Value = 57.18
$ws.Range("V31").Value = -3.1
$ws.Range("W31").Value = -1.75
$ws.Range("X31").Value = -0.74
$ws.Range("Y31").Value = -1.14
$ws.Range("V34").Value = 120
$ws.Range("W34").Value = 155
$ws.Range("X34").Value = 178.15
$ws.Range("Y34").Value = 61.98
$ws.Range("Z34").Value = 38.52
$ws.Range("V35").Value = 13.86
$ws.Range("W35").Value = 80.03
$ws.Range("X35").Value = 127.49
$ws.Range("V36").Value = 6.05
$ws.Range("W36").Value = 3.22
$ws.Range("X36").Value = 3.6
$ws.Range("Y36").Value = 3.66
$ws.Range("Z36").Value = 2.89
$ws.Range("V39").Value = 97.73999999999999
$ws.Range("W39").Value = 70.27
$ws.Range("X39").Value = 46.5
$ws.Range("Y39").Value = 57.46
$ws.Range("V41").Value = -2.35
$ws.Range("W41").Value = -1.48
$ws.Range("X41").Value = -0.5600000000000001
$ws.Range("Y41").Value = -0.86
$ws.Range("V44").Value = 100.9
$ws.Range("W44").Value = 72.69
$ws.Range("X44").Value = 47.97
$ws.Range("Y44").Value = 59.63
$ws.Range("V46").Value = 0.8100000000000001
$ws.Range("W46").Value = 0.95
$ws.Range("X46").Value = 0.91
$ws.Range("Y46").Value = 1.31
$ws.Range("V49").Value = 105.36
$ws.Range("W49").Value = 75.13
$ws.Range("X49").Value = 49.02
$ws.Range("Y49").Value = 60.62
$ws.Range("V51").Value = 5.27
$ws.Range("W51").Value = 3.38
$ws.Range("X51").Value = 1.96
$ws.Range("Y51").Value = 2.3
$ws.Range("V54").Value = 102.13
$ws.Range("W54").Value = 72.47
$ws.Range("X54").Value = 48.47
$ws.Range("Y54").Value = 60.56
$ws.Range("Z54").Value = 37.04
$ws.Range("V56").Value = 2.04
$ws.Range("W56").Value = 0.72
$ws.Range("X56").Value = 1.41
$ws.Range("Y56").Value = 2.24
$ws.Range("Z56").Value = 1.41
$ws.Range("V59").Value = 104.92
$ws.Range("W59").Value = 75.68000000000001
$ws.Range("X59").Value = 49.7
$ws.Range("Y59").Value = 61.85
$ws.Range("V61").Value = 4.83
$ws.Range("W61").Value = 3.94
$ws.Range("X61").Value = 2.63
$ws.Range("Y61").Value = 3.53
$ws.Range("V64").Value = 106.59
$ws.Range("W64").Value = 76.90000000000001
$ws.Range("X64").Value = 50.28
$ws.Range("Y64").Value = 62.44
$ws.Range("V66").Value = 6.5
$ws.Range("W66").Value = 5.15
$ws.Range("X66").Value = 3.22
$ws.Range("Y66").Value = 4.12
$ws.Range("V69").Value = 107.62
$ws.Range("W69").Value = 77.81999999999999
$ws.Range("X69").Value = 51.1
$ws.Range("Y69").Value = 63.53
$ws.Range("V71").Value = 7.53
$ws.Range("W71").Value = 6.07
$ws.Range("X71").Value = 4.04
$ws.Range("Y71").Value = 5.21
$ws.Range("V74").Value = 104.81
$ws.Range("W74").Value = 75.68000000000001
$ws.Range("X74").Value = 49.64
$ws.Range("Y74").Value = 61.71
$ws.Range("V76").Value = 4.72
$ws.Range("W76").Value = 3.94
$ws.Range("X76").Value = 2.58
$ws.Range("Y76").Value = 3.39
$ws.Range("V79").Value = 105.65
$ws.Range("W79").Value = 76.25
$ws.Range("X79").Value = 49.87
$ws.Range("Y79").Value = 61.93
$ws.Range("V81").Value = 5.56
$ws.Range("W81").Value = 4.5
$ws.Range("X81").Value = 2.81
$ws.Range("Y81").Value = 3.61
$ws.Range("V84").Value = 101.41
$ws.Range("W84").Value = 72.47
$ws.Range("X84").Value = 48.12
$ws.Range("Y84").Value = 60.5
$ws.Range("Z84").Value = 37.12
$ws.Range("V86").Value = 1.32
$ws.Range("W86").Value = 0.72
$ws.Range("X86").Value = 1.06
$ws.Range("Y86").Value = 2.18
$ws.Range("Z86").Value = 1.48
$ws.Range("V89").Value = 96.98999999999999
$ws.Range("W89").Value = 70
$ws.Range("X89").Value = 46.32
$ws.Range("Y89").Value = 57.18
$ws.Range("V91").Value = -3.1
$ws.Range("W91").Value = -1.75
$ws.Range("X91").Value = -0.74
$ws.Range("Y91").Value = -1.14
